$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.55804431438446
$ws.Range("B1").Value = 1.334033727645874
$ws.Range("C1").Value = 5.219929218292236
$ws.Range("D1").Value = 3.34108304977417
$ws.Range("E1").Value = 0.6147852540016174
